$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update total_registros values (column B) for rows 2-13
$ws.Range("B2").Value = 271
$ws.Range("B3").Value = 209
$ws.Range("B4").Value = 204
$ws.Range("B5").Value = 190
$ws.Range("B6").Value = 188
$ws.Range("B7").Value = 153
$ws.Range("B8").Value = 133
$ws.Range("B9").Value = 129
$ws.Range("B10").Value = 129
$ws.Range("B11").Value = 126
$ws.Range("B12").Value = 124
$ws.Range("B13").Value = 116

# Swap names in rows 5/6 (CHUNGA DE LA CRUZ <-> FIESTAS PERICHE)
$ws.Range("A5").Value = "FIESTAS PERICHE VIVIANA LISSETH"
$ws.Range("A6").Value = "CHUNGA DE LA CRUZ ROSA LILIANA"

# Swap names in rows 10/11 (PAIVA PINDAY <-> JIMENEZ GUERRERO)
$ws.Range("A10").Value = "JIMENEZ GUERRERO JUAN RICARDO"
$ws.Range("A11").Value = "PAIVA PINDAY ALICIA"
